$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New key/value rows to insert starting at row 4
$rows = @(
    @("options", "OPTIONS"),
    @("music", "MUSIC"),
    @("sound", "SOUND"),
    @("on", "ON"),
    @("off", "OFF"),
    @("close", "CLOSE"),
    @("yes", "YES"),
    @("no", "NO"),
    @("confirmDesc", "Are you sure?")
)

$r = 4
foreach ($pair in $rows) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r++
}

# Final row: value (RESTART) was entered before the key (confirmTitleRestart),
# so the shared string table picks up RESTART first.
$ws.Cells.Item(13, 2).Value = "RESTART"
$ws.Cells.Item(13, 1).Value = "confirmTitleRestart"

# Update the selection to match the final state of the sheet
$ws.Range("A14").Select()
